# Editados backgrounds do interior do Castillo.
# Fill in the previously-empty "Offset da Paleta" (column D) values for the
# Castillo related rows (Red Sun and Blue Moon versions), and mark the
# corresponding "Editado" (column F) cells as "Sim" where a palette offset
# was added but the row wasn't flagged as edited yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Red Sun - EletroVilla - Edifício da Jomon
$ws.Range("D7").Value = "0x50A0B8"

# Red Sun - Castillo - Gifts
$ws.Range("D10").Value = "0x529428"
$ws.Range("F10").Value = "Sim"

# Red Sun - Castillo - Restaurant / Photo Studio
$ws.Range("D11").Value = "0x5382C4"

# Blue Moon - EletroVilla - Edifício da Jomon
$ws.Range("D18").Value = "0x509D70"

# Blue Moon - Castillo - Gifts
$ws.Range("D21").Value = "0x5290D0"
$ws.Range("F21").Value = "Sim"

# Blue Moon - Castillo - Restaurant / Photo Studio
$ws.Range("D22").Value = "0x5381E8"
